$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44. This pushes the old row 44 (blank
# template row) through row 47 (sum [working weeks]) down to rows 45-48, and
# keeps all of their formulas correctly re-pointed (e.g. SUM(F2:F44) becomes
# SUM(F2:F45)).
$ws.Range("A44:G44").Insert()

# --- Row 43 (existing entry) was edited: end time changed, day changed ---
$ws.Range("C43").Value2 = 4
$ws.Range("E43").Value2 = 0.48958333333333331

# --- Row 44 is now a new data row with a fresh time entry ---
$ws.Range("A44").Value2 = 2014
$ws.Range("B44").Value2 = 3
$ws.Range("C44").Value2 = 4
$ws.Range("D44").Value2 = 0.53125
$ws.Range("E44").Value2 = 0.5625
$ws.Range("F44").Formula = "=(E44-D44)*24*60"
$ws.Range("G44").Formula = "=F44/60"

# --- Update the view state to match where the user ended up scrolled to / selecting ---
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("F44").Select()
